$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I = 9 (X), Column J = 10 (Y), Column K = 11 (locus)
# For every data row (2..131):
#   - if locus (K) is "*", Y (J) should also be "*" (mirrors X which is already "*")
#   - if locus (K) is "TR", X (I) becomes "NOM" and Y (J) becomes "ACC"
for ($r = 2; $r -le 131; $r++) {
    $locus = $ws.Cells.Item($r, 11).Value2
    if ($locus -eq "*") {
        $ws.Cells.Item($r, 10).Value2 = "*"
    } elseif ($locus -eq "TR") {
        $ws.Cells.Item($r, 9).Value2 = "NOM"
        $ws.Cells.Item($r, 10).Value2 = "ACC"
    }
}

# Reset the view: scroll back to the top-left corner and clear any lingering selection
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
